# Updates crypto price/volume data (and a few coin-name/link row shifts)
# to match the "Updated symbol list on Mon Jan 30 07:51:16 UTC 2023" commit.
#
# All target cells are plain text cells (t="inlineStr" in the original XML),
# including the numeric-looking Price/Volume columns. Assigning a bare numeric
# string to Range.Value lets Excel auto-coerce it to a Number, which would
# change the cell type and can silently drop formatting (e.g. "39.20" -> 39.2).
# Prefixing the literal with a single-quote forces Excel to store it as text
# (matching the workbook's original inlineStr cells), and resetting the style
# back to "Normal" afterwards clears the quote-prefix display flag so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'313.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.47%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.61%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.128"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.13%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08126"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.12%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'GateToken"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'4.495"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.99%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'FTXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.959"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.70%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'KuCoinToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'8.285"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.70%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'MXToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.9392"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.38%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1318"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-6.79%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.63%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.09008"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.22%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03489"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.26%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09711"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.05%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001408"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.34%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.006269"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.90%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.571"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-8.59%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.168"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-5.62%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.32%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-3.11%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'5.71%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2492"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.67%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04374"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.08%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001244"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.15%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004725"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.45%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'199.48%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-7.67%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02215"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'6.60%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05218"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.22%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007615"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.49%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01037"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.90%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1392"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.06%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-1.32%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009113"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.72%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006711"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.62%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'16.33%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = "Normal"
